$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- SaleHeader (column D/E): insert CustomerId after SaleId, shift No down,
#     and append a new Description field ---
$ws.Range("D3").Value = "CustomerId"
$ws.Range("E3").Value = "Guid"
$ws.Range("D4").Value = "No"
$ws.Range("E4").Value = "String"
$ws.Range("D5").Value = "Description"
$ws.Range("E5").Value = "String"

# --- SaleLine (column G/H): insert Discount/Tax fields before the trailing Amount ---
$ws.Range("G6").Value = "Discount"
$ws.Range("H6").Value = "decimal"
$ws.Range("G7").Value = "Tax"
$ws.Range("H7").Value = "Decimal"
$ws.Range("G8").Value = "Amount"
$ws.Range("H8").Value = "decimal"

# --- Product (column J/K): insert Barcode after ProductId, shift Name/Price down ---
$ws.Range("J3").Value = "Barcode"
$ws.Range("K3").Value = "String"
$ws.Range("J4").Value = "Name"
$ws.Range("K4").Value = "String"
$ws.Range("J5").Value = "Price"
$ws.Range("K5").Value = "Decimal"

# --- New Customer entity (column P/Q) ---
$ws.Range("P1").Value = "Customer"
$ws.Range("P1").Font.Bold = $true

$ws.Range("P2").Value = "CustomerId"
$ws.Range("Q2").Value = "Guid"

$ws.Range("P3").Value = "Name"
$ws.Range("Q3").Value = "String"

$ws.Range("P4").Value = "Surname"
$ws.Range("Q4").Value = "String"

$ws.Range("P5").Value = "Country"
$ws.Range("Q5").Value = "String"

$ws.Range("P6").Value = "City"
$ws.Range("Q6").Value = "String"

$ws.Range("P7").Value = "District"
$ws.Range("Q7").Value = "String"

$ws.Range("P8").Value = "Street"
$ws.Range("Q8").Value = "String"

# --- Column widths: D and P both need to fit the widest entry ("CustomerId") ---
$ws.Columns("D").ColumnWidth = 9.59
$ws.Columns("P").ColumnWidth = 9.59

# --- Restore the active selection to match the new working area ---
$ws.Range("I18").Select()
